$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1532.125
$ws.Range("I115").Value = 1532.125
$ws.Range("K115").Value = 4596.375
$ws.Range("M115").Value = -3029.375

$ws.Range("H116").Value = 3072.9
$ws.Range("I116").Value = 3068.4285
$ws.Range("J116").Value = 3083.3333
$ws.Range("K116").Value = 3068.4285
$ws.Range("L116").Value = 3083.3333
$ws.Range("M116").Value = 373.5715
$ws.Range("N116").Value = -9967.3333

$ws.Range("H132").Value = 3399.3235
$ws.Range("I132").Value = 3199.3333
$ws.Range("K132").Value = 9597.999899999999
$ws.Range("M132").Value = -7067.999899999999

$ws.Range("H138").Value = 1898.875
$ws.Range("J138").Value = 2066.8667
$ws.Range("L138").Value = 6200.6001
$ws.Range("N138").Value = -16480.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2435.923
$ws.Range("I61").Value = 1768.04
$ws.Range("J61").Value = 3628.5715
$ws.Range("K61").Value = 1768.04
$ws.Range("L61").Value = 3628.5715
$ws.Range("M61").Value = -1556.04
$ws.Range("N61").Value = -4052.5715

$ws.Range("H132").Value = 3182.1892
$ws.Range("I132").Value = 3001.5925
$ws.Range("K132").Value = 9004.7775
$ws.Range("M132").Value = -6474.7775

$ws.Range("H136").Value = 2435.923
$ws.Range("I136").Value = 1768.04
$ws.Range("J136").Value = 3628.5715
$ws.Range("K136").Value = 5304.12
$ws.Range("L136").Value = 10885.7145
$ws.Range("M136").Value = -2754.12
$ws.Range("N136").Value = -15985.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66668780
$ws.Range("I86").Value = 90911080
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 90911080
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -90909957
$ws.Range("N86").Value = -4696

$ws.Range("H89").Value = 66668780
$ws.Range("I89").Value = 90911080
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 454555400
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -454549784
$ws.Range("N89").Value = -23482

$ws.Range("H134").Value = 2669.848
$ws.Range("I134").Value = 2441.359
$ws.Range("J134").Value = 3942.8572
$ws.Range("K134").Value = 7324.076999999999
$ws.Range("L134").Value = 11828.5716
$ws.Range("M134").Value = -4789.076999999999
$ws.Range("N134").Value = -16898.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 34000
$ws.Range("J80").Value = 34000
$ws.Range("L80").Value = 34000
$ws.Range("N80").Value = -36246

$ws.Range("H83").Value = 34000
$ws.Range("J83").Value = 34000
$ws.Range("L83").Value = 102000
$ws.Range("N83").Value = -113232

$ws.Range("H100").Value = 49995
$ws.Range("J100").Value = 49995
$ws.Range("L100").Value = 49995
$ws.Range("N100").Value = -52159

$ws.Range("H132").Value = 4763303
$ws.Range("I132").Value = 853.37036
$ws.Range("K132").Value = 2560.11108
$ws.Range("M132").Value = -30.11108000000013

$ws.Range("H134").Value = 1615.1177
$ws.Range("I134").Value = 1449.9286
$ws.Range("K134").Value = 4349.7858
$ws.Range("M134").Value = -1814.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2872.7273
$ws.Range("J60").Value = 3120
$ws.Range("L60").Value = 9360
$ws.Range("N60").Value = -9862

$ws.Range("H113").Value = 965.04346
$ws.Range("I113").Value = 461.25
$ws.Range("J113").Value = 1233.7333
$ws.Range("K113").Value = 1383.75
$ws.Range("L113").Value = 3701.199900000001
$ws.Range("M113").Value = 786.25
$ws.Range("N113").Value = -8041.199900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1690.5
$ws.Range("I31").Value = 1690.5
$ws.Range("K31").Value = 1690.5
$ws.Range("M31").Value = -1398.5

$ws.Range("H37").Value = 1690.5
$ws.Range("I37").Value = 1690.5
$ws.Range("K37").Value = 1690.5
$ws.Range("M37").Value = -1413.5

$ws.Range("H70").Value = 5397.5
$ws.Range("I70").Value = 5353.4194
$ws.Range("J70").Value = 5456.913
$ws.Range("K70").Value = 5353.4194
$ws.Range("L70").Value = 5456.913
$ws.Range("M70").Value = -5083.4194
$ws.Range("N70").Value = -5996.913

$ws.Range("H73").Value = 5397.5
$ws.Range("I73").Value = 5353.4194
$ws.Range("J73").Value = 5456.913
$ws.Range("K73").Value = 5353.4194
$ws.Range("L73").Value = 5456.913
$ws.Range("M73").Value = -4417.4194
$ws.Range("N73").Value = -7328.913

$ws.Range("H126").Value = 3502.7
$ws.Range("I126").Value = 3002.8
$ws.Range("J126").Value = 4002.6
$ws.Range("K126").Value = 9008.400000000001
$ws.Range("L126").Value = 12007.8
$ws.Range("M126").Value = -6538.400000000001
$ws.Range("N126").Value = -16947.8

$ws.Range("H132").Value = 2170.4092
$ws.Range("I132").Value = 1882.8572
$ws.Range("J132").Value = 3288.6667
$ws.Range("K132").Value = 5648.571599999999
$ws.Range("L132").Value = 9866.000100000001
$ws.Range("M132").Value = -3118.571599999999
$ws.Range("N132").Value = -14926.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4151.875
$ws.Range("I7").Value = 2866.6667
$ws.Range("J7").Value = 4923
$ws.Range("K7").Value = 2866.6667
$ws.Range("L7").Value = 4923
$ws.Range("M7").Value = -2754.6667
$ws.Range("N7").Value = -5147

$ws.Range("H93").Value = 8029.263
$ws.Range("I93").Value = 11368.546
$ws.Range("J93").Value = 3437.75
$ws.Range("K93").Value = 11368.546
$ws.Range("L93").Value = 3437.75
$ws.Range("M93").Value = -10120.546
$ws.Range("N93").Value = -5933.75

$ws.Range("H126").Value = 4151.875
$ws.Range("I126").Value = 2866.6667
$ws.Range("J126").Value = 4923
$ws.Range("K126").Value = 8600.000100000001
$ws.Range("L126").Value = 14769
$ws.Range("M126").Value = -6130.000100000001
$ws.Range("N126").Value = -19709

$ws.Range("H132").Value = 4488.4116
$ws.Range("I132").Value = 3652.1
$ws.Range("K132").Value = 10956.3
$ws.Range("M132").Value = -8426.299999999999

$ws.Range("H140").Value = 99429
$ws.Range("J140").Value = 99429
$ws.Range("L140").Value = 99429
$ws.Range("N140").Value = -109789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1549.8667
$ws.Range("I126").Value = 1325.8
$ws.Range("K126").Value = 3977.4
$ws.Range("M126").Value = -1507.4

$ws.Range("H132").Value = 4904647
$ws.Range("I132").Value = 3748.3076
$ws.Range("K132").Value = 11244.9228
$ws.Range("M132").Value = -8714.9228

$ws.Range("H136").Value = 2169.3125
$ws.Range("I136").Value = 1919.5667
$ws.Range("J136").Value = 2585.5557
$ws.Range("K136").Value = 5758.7001
$ws.Range("L136").Value = 7756.6671
$ws.Range("M136").Value = -3208.7001
$ws.Range("N136").Value = -12856.6671

